# The source sheet has two adjacent observation records (row 2 and row 3)
# that were swapped: row 2's data moves to row 3 and vice versa. Row 4 is
# untouched. Only the columns whose row-2/row-3 values actually differ are
# written, so unrelated formatting/typing on unaffected cells is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 gets what used to be in row 3 ---
$ws.Range("A2").Value  = 130979906
$ws.Range("B2").Value  = 57884
$ws.Range("D2").Value  = "NT"
$ws.Range("E2").Value  = 100109
$ws.Range("F2").Value  = "Tretåig hackspett"
$ws.Range("G2").Value  = "Picoides tridactylus"
$ws.Range("H2").Value  = "(Linnaeus, 1758)"
$ws.Range("K2").Value  = ""
$ws.Range("L2").Value  = ""
$ws.Range("M2").Value  = "färska spår"
$ws.Range("N2").Value  = ""
$ws.Range("Q2").Value  = 591163
$ws.Range("R2").Value  = 6963104
$ws.Range("Z2").Value  = "12:38"
$ws.Range("AB2").Value = "12:38"
$ws.Range("AC2").Value = "färska ringhack på gran"

# --- Row 3 gets what used to be in row 2 ---
$ws.Range("A3").Value  = 130979911
$ws.Range("B3").Value  = 80252
$ws.Range("D3").Value  = "LC"
$ws.Range("E3").Value  = 6456
$ws.Range("F3").Value  = "Skinnlav"
$ws.Range("G3").Value  = "Leptogium saturninum"
$ws.Range("H3").Value  = "(Dicks.) Nyl."
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("Q3").Value  = 591152
$ws.Range("R3").Value  = 6963132
$ws.Range("Z3").Value  = "12:00"
$ws.Range("AB3").Value = "12:00"
$ws.Range("AC3").ClearContents()
